# Update NATMI TPM-derived LR-pair statistics on the active worksheet.
# Only the numeric columns E..T for data rows 2..7 are affected; columns
# A..D (Sending cluster, Ligand symbol, Receptor symbol, Target cluster)
# and I..L (which are unchanged in the source data) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value
$updates = @{
    2 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'M' = 2.560821
        'N' = 7.682463
        'O' = 0.2532146800817753
        'P' = 0.2532146800817753
        'Q' = 1.389092596847
        'R' = 12.501833371623
        'S' = 0.2532146800817753
        'T' = 0.2532146800817753
    }
    3 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'O' = 0.1599893367513387
        'P' = 0.1599893367513387
        'Q' = 0.8776742453635555
        'R' = 7.899068208271999
        'S' = 0.1599893367513387
        'T' = 0.1599893367513387
    }
    4 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'M' = 0.7202793333333334
        'N' = 2.160838
        'O' = 0.07122141725622931
        'P' = 0.07122141725622931
        'Q' = 0.3907085616664445
        'R' = 3.516377054998
        'S' = 0.07122141725622931
        'T' = 0.07122141725622931
    }
    5 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'M' = 1.309726333333333
        'N' = 3.929179
        'O' = 0.1295060976498071
        'P' = 0.1295060976498071
        'Q' = 0.7104483888287777
        'R' = 6.394035499458999
        'S' = 0.1295060976498071
        'T' = 0.1295060976498071
    }
    6 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'M' = 2.880229333333334
        'N' = 8.640688000000001
        'O' = 0.2847978633423207
        'P' = 0.2847978633423207
        'Q' = 1.562352559649778
        'R' = 14.061173036848
        'S' = 0.2847978633423207
        'T' = 0.2847978633423207
    }
    7 = @{
        'E' = 3
        'F' = 1
        'G' = 0.5424403333333333
        'H' = 1.627321
        'M' = 1.024174
        'N' = 3.072522
        'O' = 0.1012706049185289
        'P' = 0.1012706049185289
        'Q' = 0.5555532859513332
        'R' = 4.999979573561999
        'S' = 0.1012706049185289
        'T' = 0.1012706049185289
    }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
